$wb = $excel.ActiveWorkbook

# --- ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 2895.238
$ws.Range("J17").Value = 2895.238
$ws.Range("L17").Value = 8685.714
$ws.Range("N17").Value = -9021.714
$ws.Range("H19").Value = 1711.3846
$ws.Range("I19").Value = 818.375
$ws.Range("K19").Value = 818.375
$ws.Range("M19").Value = -643.375
$ws.Range("H70").Value = 5139
$ws.Range("I70").Value = 4205.3335
$ws.Range("J70").Value = 5979.3
$ws.Range("K70").Value = 12616.0005
$ws.Range("L70").Value = 17937.9
$ws.Range("M70").Value = -12346.0005
$ws.Range("N70").Value = -18477.9
$ws.Range("H73").Value = 5139
$ws.Range("I73").Value = 4205.3335
$ws.Range("J73").Value = 5979.3
$ws.Range("K73").Value = 12616.0005
$ws.Range("L73").Value = 17937.9
$ws.Range("M73").Value = -11680.0005
$ws.Range("N73").Value = -19809.9
$ws.Range("H86").Value = 3097.7144
$ws.Range("I86").Value = 2256.8572
$ws.Range("K86").Value = 2256.8572
$ws.Range("M86").Value = -1133.8572
$ws.Range("H89").Value = 3097.7144
$ws.Range("I89").Value = 2256.8572
$ws.Range("K89").Value = 11284.286
$ws.Range("M89").Value = -5668.286
$ws.Range("H100").Value = 7127.7144
$ws.Range("J100").Value = 9142.914000000001
$ws.Range("L100").Value = 9142.914000000001
$ws.Range("N100").Value = -10224.914
$ws.Range("H112").Value = 2538.111
$ws.Range("J112").Value = 2538.111
$ws.Range("L112").Value = 7614.333
$ws.Range("N112").Value = -9830.332999999999
$ws.Range("H132").Value = 1611.875
$ws.Range("I132").Value = 1427.8572
$ws.Range("J132").Value = 2900
$ws.Range("K132").Value = 4283.571599999999
$ws.Range("L132").Value = 8700
$ws.Range("M132").Value = -1753.571599999999
$ws.Range("N132").Value = -13760

# --- ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 13306.753
$ws.Range("I32").Value = 13348.088
$ws.Range("K32").Value = 13348.088
$ws.Range("M32").Value = -13061.088
$ws.Range("H45").Value = 53610.668
$ws.Range("I45").Value = 76166
$ws.Range("K45").Value = 76166
$ws.Range("M45").Value = -75789
$ws.Range("H102").Value = 1226.2222
$ws.Range("I102").Value = 1226.2222
$ws.Range("K102").Value = 1226.2222
$ws.Range("M102").Value = 395.7778000000001
$ws.Range("H130").Value = 222972.4
$ws.Range("J130").Value = 222972.4
$ws.Range("L130").Value = 222972.4
$ws.Range("N130").Value = -233012.4

# --- CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H59").Value = 44311.11
$ws.Range("I59").Value = 33333.332
$ws.Range("J59").Value = 49800
$ws.Range("K59").Value = 33333.332
$ws.Range("L59").Value = 49800
$ws.Range("M59").Value = -32188.332
$ws.Range("N59").Value = -52090
$ws.Range("H69").Value = 30000
$ws.Range("I69").Value = 30000
$ws.Range("K69").Value = 30000
$ws.Range("M69").Value = -29251
$ws.Range("H72").Value = 30000
$ws.Range("I72").Value = 30000
$ws.Range("K72").Value = 90000
$ws.Range("M72").Value = -86256
$ws.Range("H94").Value = 5242.8335
$ws.Range("J94").Value = 1464.7142
$ws.Range("L94").Value = 1464.7142
$ws.Range("N94").Value = -2366.7142
$ws.Range("H107").Value = 381.7143
$ws.Range("I107").Value = 336
$ws.Range("K107").Value = 336
$ws.Range("M107").Value = 1584
$ws.Range("H122").Value = 15737.538
$ws.Range("I122").Value = 26108.572
$ws.Range("K122").Value = 78325.716
$ws.Range("M122").Value = -75875.716
$ws.Range("H134").Value = 7865.4443
$ws.Range("I134").Value = 6344.846
$ws.Range("K134").Value = 19034.538
$ws.Range("M134").Value = -16499.538

# --- CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 1241.1666
$ws.Range("I5").Value = 779.8
$ws.Range("K5").Value = 2339.4
$ws.Range("M5").Value = -2227.4
$ws.Range("H21").Value = 274.5
$ws.Range("I21").Value = 49
$ws.Range("J21").Value = 500
$ws.Range("K21").Value = 147
$ws.Range("L21").Value = 1500
$ws.Range("M21").Value = 26
$ws.Range("N21").Value = -1846
$ws.Range("H107").Value = 4809.5
$ws.Range("J107").Value = 5572
$ws.Range("L107").Value = 16716
$ws.Range("N107").Value = -20556
$ws.Range("H135").Value = 1241.1666
$ws.Range("I135").Value = 779.8
$ws.Range("K135").Value = 7018.2
$ws.Range("M135").Value = -4483.2

# --- GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H58").Value = 16489
$ws.Range("I58").Value = 4998
$ws.Range("K58").Value = 4998
$ws.Range("M58").Value = -4721
$ws.Range("H62").Value = 48999.5
$ws.Range("J62").Value = 51999
$ws.Range("L62").Value = 51999
$ws.Range("N62").Value = -53371
$ws.Range("H65").Value = 48999.5
$ws.Range("J65").Value = 51999
$ws.Range("L65").Value = 155997
$ws.Range("N65").Value = -162861
$ws.Range("H107").Value = 1826.4546
$ws.Range("I107").Value = 1266.8334
$ws.Range("K107").Value = 1266.8334
$ws.Range("M107").Value = 653.1666
$ws.Range("H132").Value = 4445.587
$ws.Range("I132").Value = 3914.95
$ws.Range("K132").Value = 11744.85
$ws.Range("M132").Value = -9214.849999999999

# --- LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 3779.8
$ws.Range("I7").Value = 3699.5
$ws.Range("J7").Value = 3833.3333
$ws.Range("K7").Value = 3699.5
$ws.Range("L7").Value = 3833.3333
$ws.Range("M7").Value = -3587.5
$ws.Range("N7").Value = -4057.3333
$ws.Range("H21").Value = 26644
$ws.Range("J21").Value = 26644
$ws.Range("L21").Value = 26644
$ws.Range("N21").Value = -26992
$ws.Range("H68").Value = 3477.0908
$ws.Range("I68").Value = 2874.75
$ws.Range("K68").Value = 2874.75
$ws.Range("M68").Value = -2125.75
$ws.Range("H71").Value = 3477.0908
$ws.Range("I71").Value = 2874.75
$ws.Range("K71").Value = 14373.75
$ws.Range("M71").Value = -10629.75
$ws.Range("H126").Value = 3779.8
$ws.Range("I126").Value = 3699.5
$ws.Range("J126").Value = 3833.3333
$ws.Range("K126").Value = 11098.5
$ws.Range("L126").Value = 11499.9999
$ws.Range("M126").Value = -8628.5
$ws.Range("N126").Value = -16439.9999
$ws.Range("H132").Value = 6952.645
$ws.Range("I132").Value = 4890
$ws.Range("K132").Value = 14670
$ws.Range("M132").Value = -12140

# --- WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H54").Value = 130000
$ws.Range("J54").Value = 60000
$ws.Range("L54").Value = 60000
$ws.Range("N54").Value = -61040
$ws.Range("H62").Value = 29612.8
$ws.Range("I62").Value = 29612.8
$ws.Range("J62").Value = 0
$ws.Range("K62").Value = 29612.8
$ws.Range("L62").Value = 0
$ws.Range("M62").Value = -28988.8
$ws.Range("H65").Value = 29612.8
$ws.Range("I65").Value = 29612.8
$ws.Range("J65").Value = 0
$ws.Range("K65").Value = 148064
$ws.Range("L65").Value = 0
$ws.Range("M65").Value = -144944
$ws.Range("H97").Value = 67500
$ws.Range("J97").Value = 67500
$ws.Range("L97").Value = 67500
$ws.Range("N97").Value = -69482
$ws.Range("H107").Value = 4140.3184
$ws.Range("I107").Value = 4365
$ws.Range("J107").Value = 3541.1667
$ws.Range("K107").Value = 13095
$ws.Range("L107").Value = 10623.5001
$ws.Range("M107").Value = -11175
$ws.Range("N62").ClearContents()
$ws.Range("N65").ClearContents()
